$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; everything from row 31 downward shifts to row+1.
$ws.Rows.Item(31).Insert()

# Preserve the date-column number format used throughout column D.
$dateFmt = $ws.Cells.Item(32, 4).NumberFormat

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 44571
$ws.Cells.Item(31, 4).NumberFormat = $dateFmt
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100106
$ws.Cells.Item(31, 8).Value = "Oleaginosos"
$ws.Cells.Item(31, 9).Value = 100106002
$ws.Cells.Item(31, 10).Value = "Palta"
$ws.Cells.Item(31, 11).Value = "Fuerte"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 300
$ws.Cells.Item(31, 14).Value = 34000
$ws.Cells.Item(31, 15).Value = 35000
$ws.Cells.Item(31, 16).Value = 34500
$ws.Cells.Item(31, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(31, 18).Value = "Perú"
$ws.Cells.Item(31, 19).Value = 3450
$ws.Cells.Item(31, 20).Value = 10
